# Sketch a mod matrix on Sheet1.
#
# Before: two side-by-side stub tables ("Noisy"/"Wide") sharing the same
# Parm1/Parm2/Parm3 column headers, rows 2-3 only holding the LFO/Expr.
# row labels (no data).
#
# After: a single small table -
#   A1            "Noisy"
#   A2  B2    C2      D2    "" "Freq" "Cutoff" "LFOFreq"
#   A3  B3 C3 D3       "Expr."  100    500       5
#   A4  B4 C4 D4       "LFO"    10     50        0
#
# The old "Wide" side table (F:I) is removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- clear the old Parm1/Parm2/Parm3 headers and the side-by-side "Wide" table (columns F:I, rows 1-3) ---
$ws.Range("B1:D1").Clear()
$ws.Range("F1:I3").Clear()

# --- row 1: single remaining header label ---
$ws.Range("A1").Value = "Noisy"

# --- row 2: new column headers (bold, like the rest of the labels) ---
$ws.Range("A2").Value = "Scale.Coeff"
$ws.Range("B2").Value = "Freq"
$ws.Range("C2").Value = "Cutoff"
$ws.Range("D2").Value = "LFOFreq"
$ws.Range("A2:D2").Font.Bold = $true

# --- row 3: "Expr." data row ---
$ws.Range("A3").Value = "Expr."
$ws.Range("B3").Value = 100
$ws.Range("C3").Value = 500
$ws.Range("D3").Value = 5

# --- row 4: new "LFO" data row ---
$ws.Range("A4").Value = "LFO"
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 50
$ws.Range("D4").Value = 0

# Row 4 labels should be bold like the other row labels/headers.
$ws.Range("A4").Font.Bold = $true

# --- cosmetic bits from the diff: column width, zoom, selection ---
# (the stored OOXML <col width> includes ~0.8333 chars of padding over the
# COM ColumnWidth value, so 17.1667 here serializes out to width="18")
$ws.Columns.Item(1).ColumnWidth = 17.1667
$excel.ActiveWindow.Zoom = 130
$ws.Range("C5").Select()
